$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Improve trial shuffle | Added SubNum to traj output file"
#
# The traj output header row (row 1/2 on Sheet1) gets a new leading
# "sub_num" / "subject number" column inserted before the existing
# iBlock/iTrial/cat_block columns, which shift one column to the right
# (BB->BD->BE). The block-number column (iBlock, BC) stays put.
#
# Using Insert() on the BD column keeps its custom width (20.5) attached
# to the cell content ("cat_block") as it shifts right to BE, exactly like
# it would in Excel when a column is inserted before it.

# Step 1: insert a blank column at BD (56) - this pushes the old BD
# ("cat_block", with its custom width) to BE, and the old BE ("sub_num")
# to BF.
[void]$ws.Columns.Item(56).Insert()

# Step 2: the old BB content ("iTrial" / "trial number") moves into the
# now-empty BD.
$ws.Range("BD1:BD2").Value2 = $ws.Range("BB1:BB2").Value2

# Step 3: the "sub_num" / "subject number" column (now sitting in BF after
# the insert) becomes the new BB.
$ws.Range("BB1:BB2").Value2 = $ws.Range("BF1:BF2").Value2

# Step 4: drop the now-duplicate trailing column (BF).
[void]$ws.Columns.Item(58).Delete()

[void]$ws.Range("BD7").Select()
